# Update "想去人数" (column F) values across the four worksheets of the
# 广州-漫展信息 workbook, as captured by the source diff.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (cell address, new value) pairs to update.
$updates = @{
    "展览" = @(
        @{ Cell = "F2";  Value = 1914 },
        @{ Cell = "F3";  Value = 1535 },
        @{ Cell = "F4";  Value = 901 },
        @{ Cell = "F5";  Value = 798 },
        @{ Cell = "F6";  Value = 13451 },
        @{ Cell = "F7";  Value = 13294 },
        @{ Cell = "F14"; Value = 2 },
        @{ Cell = "F15"; Value = 696 },
        @{ Cell = "F17"; Value = 19 },
        @{ Cell = "F18"; Value = 68 },
        @{ Cell = "F20"; Value = 83 },
        @{ Cell = "F22"; Value = 407 },
        @{ Cell = "F23"; Value = 292 },
        @{ Cell = "F24"; Value = 296 },
        @{ Cell = "F25"; Value = 442 },
        @{ Cell = "F26"; Value = 773 },
        @{ Cell = "F27"; Value = 33 }
    )
    "演出" = @(
        @{ Cell = "F5";  Value = 136 },
        @{ Cell = "F7";  Value = 133 },
        @{ Cell = "F8";  Value = 525 },
        @{ Cell = "F11"; Value = 39 }
    )
    "本地生活" = @(
        @{ Cell = "F2"; Value = 201 },
        @{ Cell = "F3"; Value = 62 }
    )
    "全部类型" = @(
        @{ Cell = "F2";  Value = 201 },
        @{ Cell = "F3";  Value = 1914 },
        @{ Cell = "F4";  Value = 1535 },
        @{ Cell = "F5";  Value = 901 },
        @{ Cell = "F7";  Value = 798 },
        @{ Cell = "F8";  Value = 13451 },
        @{ Cell = "F9";  Value = 13294 },
        @{ Cell = "F16"; Value = 2 },
        @{ Cell = "F17"; Value = 696 },
        @{ Cell = "F21"; Value = 19 },
        @{ Cell = "F22"; Value = 68 },
        @{ Cell = "F24"; Value = 83 },
        @{ Cell = "F25"; Value = 136 },
        @{ Cell = "F28"; Value = 62 },
        @{ Cell = "F29"; Value = 407 },
        @{ Cell = "F30"; Value = 292 },
        @{ Cell = "F31"; Value = 296 },
        @{ Cell = "F32"; Value = 442 },
        @{ Cell = "F33"; Value = 773 },
        @{ Cell = "F34"; Value = 133 },
        @{ Cell = "F35"; Value = 525 },
        @{ Cell = "F38"; Value = 33 },
        @{ Cell = "F39"; Value = 39 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Range($entry.Cell).Value = $entry.Value
    }
}
